$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 8-13 on each sheet refer to the same six handed-off files.
# A fresh "Generate Report for Handoff" run updates the handoff timestamps
# and marks the handoff priority as "ht" for those rows.

$rows = 8..13

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G
    $overview.Range("G$r").Value = "2016-08-18 12:21:18"

    # zh-cn sheet: "Latest Handoff Datetime" column H, "Priority" column E
    $zhcn.Range("H$r").Value = "2016-08-18 12:21:11"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: "Latest Handoff Datetime" column H, "Priority" column E
    $dede.Range("H$r").Value = "2016-08-18 12:21:18"
    $dede.Range("E$r").Value = "ht"
}
